$d = $word.ActiveDocument

$replacements = @(
    @{old="84×39="; new="57×99="},
    @{old="42×25="; new="75×55="},
    @{old="14×48="; new="41×43="},
    @{old="20×11="; new="32×69="},
    @{old="17×32="; new="26×38="},
    @{old="82×91="; new="40×49="},
    @{old="48×52="; new="21×51="},
    @{old="86×78="; new="37×42="},
    @{old="60×76="; new="21×83="},
    @{old="26×66="; new="16×30="},
    @{old="59×51="; new="51×94="},
    @{old="60×17="; new="86×16="},
    @{old="71×95="; new="26×22="},
    @{old="88×33="; new="46×99="},
    @{old="46×95="; new="44×43="},
    @{old="46×67="; new="23×26="},
    @{old="88×72="; new="54×89="},
    @{old="85×37="; new="68×67="},
    @{old="32×60="; new="75×25="},
    @{old="28×84="; new="68×67="},
    @{old="77×46="; new="46×39="},
    @{old="17×35="; new="16×11="},
    @{old="94×92="; new="27×22="},
    @{old="44×76="; new="18×68="},
    @{old="54×51="; new="62×98="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
